$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (C, D) for "Employee Age" and "Gender(M/F)",
# pushing the existing columns (starting with "Date of Joining") to the right.
$ws.Range("C1:D1").EntireColumn.Insert()

# The original sheet had "Date of Joining" twice (columns C and J).
# Column C's occurrence survives (now shifted to E); remove the duplicate,
# which is now located in column L.
$ws.Range("L1").EntireColumn.Delete()

# Remove the old "Gross Salary" header; it is being replaced by a renamed /
# relocated "Net Salary" column plus a new "CTC" column. It is currently in
# column N.
$ws.Range("N1").EntireColumn.Delete()

# Insert a new column for "CTC" right before the remaining "Net Salary"
# column (now in column N), pushing "Net Salary" to column O.
$ws.Range("N1").EntireColumn.Insert()

# Set the text for the newly inserted / renamed header cells.
$ws.Range("C1").Value = "Employee Age"
$ws.Range("D1").Value = "Gender(M/F)"
$ws.Range("N1").Value = "CTC"
$ws.Range("O1").Value = "Net Salary (to be disbursed to the employee)"

# Match the selection state left behind in the saved workbook.
$ws.Columns.Item(4).Select()
